# Update odds values in Sheet1 to match the 2024-11-21 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Independiente - Gimnasia L.P.)
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5
$ws.Range("AX3").Value = 29

# Row 4 (Talleres Cordoba - Sarmiento Junin)
$ws.Range("G4").Value = 1.5
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("AB4").Value = 34
$ws.Range("AF4").Value = 81
$ws.Range("AH4").Value = 34
$ws.Range("AO4").Value = 7.5
$ws.Range("AT4").Value = 2.63

# Row 6 (Always Ready - Royal Pari)
$ws.Range("H6").Value = 5.5
$ws.Range("J6").Value = 1.67
$ws.Range("K6").Value = 2.6
$ws.Range("N6").Value = 12
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 4.33
$ws.Range("Q6").Value = 1.67
$ws.Range("R6").Value = 2.15
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 6
$ws.Range("Y6").Value = 10
$ws.Range("AC6").Value = 12
$ws.Range("AE6").Value = 26
$ws.Range("AF6").Value = 81
$ws.Range("AG6").Value = 26
$ws.Range("AK6").Value = 81
$ws.Range("AN6").Value = 3.1
$ws.Range("AP6").Value = 19
$ws.Range("AZ6").Value = 301
$ws.Range("BA6").Value = 301

# Row 7 (Tomayapo - Bolivar)
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 15
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 2.1

# Row 9 (Vasco - Internacional)
$ws.Range("G9").Value = 3.5
$ws.Range("L9").Value = 2.88
$ws.Range("X9").Value = 17
$ws.Range("AW9").Value = 4

# Row 10 (Once Caldas - Deportes Tolima)
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("Q10").Value = 2.7
$ws.Range("R10").Value = 1.44

# Row 13 (Sportivo Trinidense - Cerro Porteno)
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("O13").Value = 1.33
$ws.Range("P13").Value = 3.25
